$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: PE Reval - switch from "linear vs D2" formula to a chained
#     revaluation formula that blends period-over-period change with a
#     sensitivity factor in $B3 ---
$ws.Range("D3").Formula = '=MIN(1,C3*(1+(D$2/C$2-1)*$B3))'
$ws.Range("E3").Formula = '=MIN(1,D3*(1+(E$2/D$2-1)*$B3))'
$ws.Range("F3:AB3").Formula = '=MIN(1,E3*(1+(F$2/E$2-1)*$B3))'

# --- Row 4: Credit Reval - same chained-formula change ---
$ws.Range("D4").Formula = '=MIN(1,C4*(1+(D$2/C$2-1)*$B4))'
$ws.Range("E4").Formula = '=MIN(1,D4*(1+(E$2/D$2-1)*$B4))'
$ws.Range("F4:AB4").Formula = '=MIN(1,E4*(1+(F$2/E$2-1)*$B4))'

# --- Row 5: RE Reval - same chained-formula change ---
$ws.Range("D5").Formula = '=MIN(1,C5*(1+(D$2/C$2-1)*$B5))'
$ws.Range("E5").Formula = '=MIN(1,D5*(1+(E$2/D$2-1)*$B5))'
$ws.Range("F5:AB5").Formula = '=MIN(1,E5*(1+(F$2/E$2-1)*$B5))'

# --- Row 60: Cash Flow according to PM - Credit: updated estimates ---
$ws.Range("D60:F60").Value = 60
$ws.Range("G60:Q60").Value = 25
$ws.Range("R60:AB60").Value = -12

# --- Row 65: Credit Return Assumed (Annual): new estimate (3/30 vintage) ---
$ws.Range("D65:O65").Value = 0.08
$ws.Range("P65:AB65").Value = 0.1

# --- Move the selection/scroll position back up to A8 (was parked at A66
#     with the view scrolled down to A36) ---
$ws.Range("A8").Select()
